$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows (at what become rows 28, 30, 31) to make room
# for "sesenta", "setenta", "ochenta" and their Quechua equivalents.
$ws.Rows("28").Insert()
$ws.Rows("30").Insert()
$ws.Rows("31").Insert()

# New Spanish values (column A) entered first
$ws.Range("A28").Value = "sesenta"
$ws.Range("A30").Value = "setenta"
$ws.Range("A31").Value = "ochenta"

# Append three new rows at the bottom for quinientos / mil / un millón
$ws.Range("A40").Value = "quinientos"
$ws.Range("A41").Value = "mil"
$ws.Range("A42").Value = "un millón"

# New Quechua values (column B) entered afterwards
$ws.Range("B28").Value = "suqta chunka"
$ws.Range("B30").Value = "qanchis chunka"
$ws.Range("B31").Value = "pusaq chunka"

# Correct B9 ("ocho" -> "pusaq" instead of duplicated "qanchis")
$ws.Range("B9").Value = "pusaq"

# Clear the cell/column formatting (removes the border style and the
# explicit column style/width that were set before).
$ws.Cells.ClearFormats()

# Restore the view/selection state
$ws.Range("A4").Select()
$ws.Range("J11").Select()
